# Auto-generated Excel COM-interop script to apply cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force "Price" (column D) cells whose new value looks like a plain decimal number to remain
# stored as TEXT (matching the original workbook, where every Price cell is an inline string),
# instead of being auto-converted to a numeric value by Excel.
$textPriceCells = 'D4', 'D5', 'D6', 'D7', 'D8', 'D9', 'D10', 'D11', 'D15', 'D17', 'D18', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D26', 'D27', 'D28', 'D29', 'D31', 'D33', 'D37', 'D38', 'D39', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D48', 'D49', 'D51'
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Per-row "Price" (D) and "Volume(1h)" (E) updates ---
$ws.Range("D2").Value = '27.719.02'
$ws.Range("E2").Value = '  +1.96%  '
$ws.Range("D3").Value = '1.561.04'
$ws.Range("E3").Value = '  -0.74%  '
$ws.Range("D4").Value = '0.987'
$ws.Range("E4").Value = '  -2.69%  '
$ws.Range("D5").Value = '210.38'
$ws.Range("E5").Value = '  -0.71%  '
$ws.Range("D6").Value = '0.494'
$ws.Range("E6").Value = '  +0.28%  '
$ws.Range("D7").Value = '0.985'
$ws.Range("E7").Value = '  -2.62%  '
$ws.Range("D8").Value = '23.28'
$ws.Range("E8").Value = '  +5.41%  '
$ws.Range("D9").Value = '0.249'
$ws.Range("E9").Value = '  -0.28%  '
$ws.Range("D10").Value = '0.0593'
$ws.Range("E10").Value = '  -0.91%  '
$ws.Range("D11").Value = '0.0873'
$ws.Range("E11").Value = '  +0.86%  '
$ws.Range("D12").Value = '1.783.02'
$ws.Range("E12").Value = '  -0.71%  '
$ws.Range("D13").Value = '1.570.30'
$ws.Range("E13").Value = '  +0.98%  '
$ws.Range("E14").Value = '  -1.26%  '
$ws.Range("D15").Value = '0.518'
$ws.Range("E15").Value = '  -0.51%  '
$ws.Range("D16").Value = '27.681.85'
$ws.Range("E16").Value = '  +1.85%  '
$ws.Range("D17").Value = '63.01'
$ws.Range("D18").Value = '229.17'
$ws.Range("E18").Value = '  +6.08%  '
$ws.Range("D21").Value = '0.990'
$ws.Range("E21").Value = '  -2.31%  '
$ws.Range("D22").Value = '4.09'
$ws.Range("E22").Value = '  -1.35%  '
$ws.Range("D23").Value = '9.32'
$ws.Range("E23").Value = '  +1.14%  '
$ws.Range("D24").Value = '1.93'
$ws.Range("E24").Value = '  -1.45%  '
$ws.Range("D25").Value = '150.14'
$ws.Range("E25").Value = '  -2.95%  '
$ws.Range("D26").Value = '15.23'
$ws.Range("E26").Value = '  +0.66%  '
$ws.Range("D27").Value = '6.56'
$ws.Range("E27").Value = '  -0.77%  '
$ws.Range("D28").Value = '0.107'
$ws.Range("E28").Value = '  +1.00%  '
$ws.Range("D29").Value = '0.987'
$ws.Range("E29").Value = '  -2.63%  '
$ws.Range("E30").Value = '  -1.07%  '
$ws.Range("D31").Value = '0.0470'
$ws.Range("E31").Value = '  -0.85%  '
$ws.Range("E32").Value = '  -0.74%  '
$ws.Range("D33").Value = '3.12'
$ws.Range("E33").Value = '  -2.25%  '
$ws.Range("D34").Value = '1.406.72'
$ws.Range("E34").Value = '  -1.99%  '
$ws.Range("E35").Value = '  -2.36%  '
$ws.Range("E36").Value = '  -4.87%  '
$ws.Range("D37").Value = '2.30'
$ws.Range("E37").Value = '  -2.72%  '
$ws.Range("D38").Value = '0.0167'
$ws.Range("E38").Value = '  -0.21%  '
$ws.Range("D39").Value = '0.541'
$ws.Range("E39").Value = '  +1.85%  '
$ws.Range("E40").Value = '  +2.61%  '
$ws.Range("D41").Value = '0.807'
$ws.Range("E41").Value = '  -0.41%  '
$ws.Range("D42").Value = '5.63'
$ws.Range("E42").Value = '  -3.79%  '
$ws.Range("D43").Value = '0.988'
$ws.Range("E43").Value = '  -2.72%  '
$ws.Range("D44").Value = '1.84'
$ws.Range("E44").Value = '  +5.29%  '
$ws.Range("D45").Value = '0.966'
$ws.Range("E45").Value = '  -4.01%  '
$ws.Range("D46").Value = '63.77'
$ws.Range("E46").Value = '  -1.52%  '
$ws.Range("D47").Value = '1.693.37'
$ws.Range("E47").Value = '  -0.97%  '
$ws.Range("D48").Value = '86.38'
$ws.Range("D49").Value = '0.0522'
$ws.Range("E49").Value = '  +1.00%  '
$ws.Range("D50").Value = '0.0₇0990'
$ws.Range("E50").Value = '  -2.41%  '

# --- Rows 19 & 20 swapped their coin (Chainlink <-> ShibaInu) along with price/volume ---
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '0.0₃0704'
$ws.Range("E19").Value = '  -0.12%  '

$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").Value = '7.43'
$ws.Range("E20").Value = '  +0.45%  '

# --- Row 51: Algorand replaced by BitcoinSV ---
$ws.Range("B51").Value = 'BitcoinSV'
$ws.Range("C51").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D51").Value = '39.93'
$ws.Range("E51").Value = '  +17.20%  '
